# Fruta / hortaliza, semanal
#
# A new weekly 3-row block (date 44540) is inserted into the daily price
# table right before the existing 2021-09-10 block (old row 60), pushing
# every subsequent row down by three positions (old A1:T144 -> A1:T147).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 60; everything from old row 60 downward
# shifts down to row 63 onward (dimension grows from T144 to T147).
$ws.Range("60:62").EntireRow.Insert()

# Fill in the new block (rows 60-62) with the new week's data. Columns
# A, B, C, E, F, G, H, I, J, K are constant across the whole sheet, so
# Excel's insert already carried usable blanks/styles for the row; we
# only need to populate them explicitly since the inserted rows are
# otherwise empty.

# Row 60: Especial
$ws.Range("A60").Value = 8
$ws.Range("B60").Value = "Terminal La Palmera de La Serena"
$ws.Range("C60").Value = "Coquimbo"
$ws.Range("D60").Value = 44540
$ws.Range("E60").Value = 4
$ws.Range("F60").Value = "Fruta"
$ws.Range("G60").Value = 100107
$ws.Range("H60").Value = "Otros"
$ws.Range("I60").Value = 100107002
$ws.Range("J60").Value = "Chirimoya"
$ws.Range("K60").Value = "Cultivar IV Región"
$ws.Range("L60").Value = "Especial"
$ws.Range("M60").Value = 240
$ws.Range("N60").Value = 1600
$ws.Range("O60").Value = 1700
$ws.Range("P60").Value = 1650
$ws.Range("Q60").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R60").Value = "Provincia de Limarí"
$ws.Range("S60").Value = 1650
$ws.Range("T60").Value = 1

# Row 61: Primera
$ws.Range("A61").Value = 8
$ws.Range("B61").Value = "Terminal La Palmera de La Serena"
$ws.Range("C61").Value = "Coquimbo"
$ws.Range("D61").Value = 44540
$ws.Range("E61").Value = 4
$ws.Range("F61").Value = "Fruta"
$ws.Range("G61").Value = 100107
$ws.Range("H61").Value = "Otros"
$ws.Range("I61").Value = 100107002
$ws.Range("J61").Value = "Chirimoya"
$ws.Range("K61").Value = "Cultivar IV Región"
$ws.Range("L61").Value = "Primera"
$ws.Range("M61").Value = 300
$ws.Range("N61").Value = 1400
$ws.Range("O61").Value = 1500
$ws.Range("P61").Value = 1450
$ws.Range("Q61").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R61").Value = "Provincia de Limarí"
$ws.Range("S61").Value = 1450
$ws.Range("T61").Value = 1

# Row 62: Segunda
$ws.Range("A62").Value = 8
$ws.Range("B62").Value = "Terminal La Palmera de La Serena"
$ws.Range("C62").Value = "Coquimbo"
$ws.Range("D62").Value = 44540
$ws.Range("E62").Value = 4
$ws.Range("F62").Value = "Fruta"
$ws.Range("G62").Value = 100107
$ws.Range("H62").Value = "Otros"
$ws.Range("I62").Value = 100107002
$ws.Range("J62").Value = "Chirimoya"
$ws.Range("K62").Value = "Cultivar IV Región"
$ws.Range("L62").Value = "Segunda"
$ws.Range("M62").Value = 240
$ws.Range("N62").Value = 900
$ws.Range("O62").Value = 1000
$ws.Range("P62").Value = 950
$ws.Range("Q62").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R62").Value = "Provincia de Limarí"
$ws.Range("S62").Value = 950
$ws.Range("T62").Value = 1
